$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '30.700.72'
$ws.Range("E2").Value = '  -0.89%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.890.55'
$ws.Range("E3").Value = '  -1.06%  '
$ws.Range("E4").Value = '  +0.13%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '236.27'
$ws.Range("E5").Value = '  -3.92%  '
$ws.Range("E6").Value = '  +0.08%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.4875'
$ws.Range("E7").Value = '  -2.41%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.2900'
$ws.Range("E8").Value = '  -3.19%  '
$ws.Range("E9").Value = '  -2.99%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '1.887.33'
$ws.Range("E10").Value = '  -1.23%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '16.72'
$ws.Range("E11").Value = '  -1.84%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.07250'
$ws.Range("E12").Value = '  -0.74%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '89.44'
$ws.Range("E13").Value = '  -2.23%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '5.011'
$ws.Range("E14").Value = '  -1.69%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.6643'
$ws.Range("E15").Value = '  -2.59%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '30.650.18'
$ws.Range("E16").Value = '  -0.95%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.000007879'
$ws.Range("E17").Value = '  -2.02%  '
$ws.Range("E18").Value = '  +0.03%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '13.00'
$ws.Range("E19").Value = '  -2.05%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '2.130.45'
$ws.Range("E20").Value = '  -1.26%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '1.001'
$ws.Range("E21").Value = '  +0.02%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '4.743'
$ws.Range("E22").Value = '  -2.89%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '191.55'
$ws.Range("E23").Value = '  +5.11%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '6.070'
$ws.Range("E24").Value = '  -1.01%  '
$ws.Range("E25").Value = '  -0.84%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '159.55'
$ws.Range("E26").Value = '  +3.06%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '18.28'
$ws.Range("E27").Value = '  -4.20%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '1.825'
$ws.Range("E28").Value = '  -6.19%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '1.403'
$ws.Range("E29").Value = '  +0.74%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '4.268'
$ws.Range("E30").Value = '  -1.90%  '
$ws.Range("E31").Value = '  +0.64%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '3.945'
$ws.Range("E32").Value = '  -2.77%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.05203'
$ws.Range("E33").Value = '  -1.07%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.7342'
$ws.Range("E34").Value = '  -2.21%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.085'
$ws.Range("E35").Value = '  -4.79%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '2.692'
$ws.Range("E36").Value = '  +0.87%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.01823'
$ws.Range("E37").Value = '  -7.78%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '2.674'
$ws.Range("E38").Value = '  -2.31%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.9235'
$ws.Range("E39").Value = '  -1.61%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '2.054'
$ws.Range("E40").Value = '  -5.77%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.4444'
$ws.Range("E41").Value = '  +1.30%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '104.68'
$ws.Range("E42").Value = '  -1.75%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '1.000'
$ws.Range("E43").Value = '  -0.04%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '5.711'
$ws.Range("E44").Value = '  -2.88%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.1336'
$ws.Range("E45").Value = '  -0.71%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '7.317'
$ws.Range("E46").Value = '  -6.39%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.4191'
$ws.Range("E47").Value = '  +7.04%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.05835'
$ws.Range("E48").Value = '  -0.25%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '8.607'
$ws.Range("E49").Value = '  -0.10%  '
$ws.Range("E50").Value = '  +0.90%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '33.22'
$ws.Range("E51").Value = '  -0.22%  '
